$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 389, shifting the existing rows 389-413 down to 390-414.
$ws.Rows(389).Insert()

# Populate the newly inserted row 389 with its values. Most columns duplicate
# what is now row 390 (the former row 389), except D, K, L, M and P which
# carry new figures for this weekly entry.
$ws.Range("A389").Value = 10
$ws.Range("B389").Value = "Vega Modelo de Temuco"
$ws.Range("C389").Value = "La Araucanía"
$ws.Range("D389").Value = 44585
$ws.Range("E389").Value = 9
$ws.Range("F389").Value = 100112043
$ws.Range("G389").Value = "Pepino ensalada"
$ws.Range("H389").Value = "Sin especificar"
$ws.Range("I389").Value = "Primera"
$ws.Range("J389").Value = 235
$ws.Range("K389").Value = 10000
$ws.Range("L389").Value = 12000
$ws.Range("M389").Value = 11064
$ws.Range("N389").Value = "$/caja 60 unidades"
$ws.Range("O389").Value = "Región del Maule"
$ws.Range("P389").Value = 184
$ws.Range("Q389").Value = 60
$ws.Range("R389").Value = "Hortaliza"
